$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "726×5=3630" "664×9=5976"
Replace-Text "511×5=2555" "284×4=1136"
Replace-Text "611×4=2444" "934×2=1868"
Replace-Text "499×9=4491" "704×7=4928"
Replace-Text "262×2=524" "965×8=7720"
Replace-Text "199×4=796" "919×7=6433"
Replace-Text "166×2=332" "708×2=1416"
Replace-Text "980×6=5880" "262×9=2358"
Replace-Text "231×6=1386" "723×9=6507"
Replace-Text "704×4=2816" "936×5=4680"
Replace-Text "798×8=6384" "161×9=1449"
Replace-Text "471×9=4239" "821×3=2463"
Replace-Text "876×2=1752" "886×6=5316"
Replace-Text "511×6=3066" "498×9=4482"
Replace-Text "713×6=4278" "772×4=3088"
Replace-Text "738×9=6642" "761×9=6849"
Replace-Text "825×2=1650" "504×4=2016"
Replace-Text "317×4=1268" "176×7=1232"
Replace-Text "888×8=7104" "472×3=1416"
Replace-Text "372×7=2604" "844×6=5064"
Replace-Text "151×6=906" "690×9=6210"
Replace-Text "812×3=2436" "124×8=992"
Replace-Text "858×6=5148" "223×3=669"
Replace-Text "621×6=3726" "291×6=1746"
Replace-Text "415×2=830" "830×4=3320"
